# Score Tracker schedule.xlsx update
# - Purse/Memorial Tournament label correction (capitalize "The")
# - Corrected start date for row 3 (Valspar Championship)
# - Selection moved to F8 as last edited cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Memorial Tournament event name capitalization (row 12, column B)
$ws.Range("B12").Value = "The Memorial Tournament "

# Correct the StartDate for row 3 (Valspar Championship)
$ws.Range("A3").Value = 45738

# Update the active selection to reflect the cell last worked on
$ws.Range("F8").Select()
